$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04941833333333333
$ws.Range("H2").Value = 0.148255
$ws.Range("I2").Value = 0.005167549122999764
$ws.Range("J2").Value = 0.005167549122999764
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 0.01434067320444445
$ws.Range("R2").Value = 0.12906605884
$ws.Range("S2").Value = 0.000177215391068647
$ws.Range("T2").Value = 0.000177215391068647
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04941833333333333
$ws.Range("H3").Value = 0.148255
$ws.Range("I3").Value = 0.005167549122999764
$ws.Range("J3").Value = 0.005167549122999764
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 0.3513816793622222
$ws.Range("R3").Value = 3.16243511426
$ws.Range("S3").Value = 0.004342211891645047
$ws.Range("T3").Value = 0.004342211891645047
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04941833333333333
$ws.Range("H4").Value = 0.148255
$ws.Range("I4").Value = 0.005167549122999764
$ws.Range("J4").Value = 0.005167549122999764
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 0.05244749596611112
$ws.Range("R4").Value = 0.472027463695
$ws.Range("S4").Value = 0.0006481218402860702
$ws.Range("T4").Value = 0.0006481218402860701
$ws.Range("I5").Value = 0.806706161560336
$ws.Range("J5").Value = 0.806706161560336
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 2.23872268256889
$ws.Range("R5").Value = 20.14850414312
$ws.Range("S5").Value = 0.0276650970306429
$ws.Range("T5").Value = 0.0276650970306429
$ws.Range("I6").Value = 0.806706161560336
$ws.Range("J6").Value = 0.806706161560336
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.6778627555178796
$ws.Range("T6").Value = 0.6778627555178796
$ws.Range("I7").Value = 0.806706161560336
$ws.Range("J7").Value = 0.806706161560336
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("Q7").Value = 8.187579285112223
$ws.Range("S7").Value = 0.1011783090118137
$ws.Range("T7").Value = 0.1011783090118137
$ws.Range("I8").Value = 0.1881262893166642
$ws.Range("J8").Value = 0.1881262893166643
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 0.5220768244364444
$ws.Range("R8").Value = 4.698691419928
$ws.Range("S8").Value = 0.006451583359539088
$ws.Range("T8").Value = 0.006451583359539089
$ws.Range("I9").Value = 0.1881262893166642
$ws.Range("J9").Value = 0.1881262893166643
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("S9").Value = 0.1580796217235908
$ws.Range("T9").Value = 0.1580796217235908
$ws.Range("I10").Value = 0.1881262893166642
$ws.Range("J10").Value = 0.1881262893166643
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("S10").Value = 0.02359508423353437
$ws.Range("T10").Value = 0.02359508423353437
